# Weekly fruit/vegetable price update: a new record for this week is
# inserted at the top of the data block (row 77, right after the rows
# that are not part of this rotating weekly block), pushing every
# existing data row down by one. The row that used to be last (174)
# becomes row 175.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 77 - shifts rows 77:174 down to 78:175
# and carries formatting down with them (same behavior as right-clicking
# the row header and choosing "Insert").
$ws.Rows(77).Insert()

# Populate the newly inserted row 77 with this week's record.
$ws.Range("A77").Value = 4
$ws.Range("B77").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C77").Value = "Los Lagos"
$ws.Range("D77").Value = 44482
$ws.Range("E77").Value = 10
$ws.Range("F77").Value = 100112040
$ws.Range("G77").Value = "Cilantro"
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 70
$ws.Range("K77").Value = 10000
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = 10000
$ws.Range("N77").Value = "$/caja 36 atados"
$ws.Range("O77").Value = "Región Metropolitana"
$ws.Range("P77").Value = 278
$ws.Range("Q77").Value = 36
$ws.Range("R77").Value = "Hortaliza"
